$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 10 new rows above row 2, pushing the existing 4 data rows down to rows 12-15
$ws.Range("A2:A11").EntireRow.Insert()

# Fill the 10 new rows (NO 1-10) with the new records
# Row 2 (NO=1)
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).NumberFormat = "@"
$ws.Cells.Item(2, 2).Value = "2022-11-16"
$ws.Cells.Item(2, 3).Value = "idong"
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "617547260712000"
$ws.Cells.Item(2, 5).Value = "Orang Pribadi"
$ws.Cells.Item(2, 6).Value = "jahanjang"
$ws.Cells.Item(2, 7).Value = "karyawan swasta"
$ws.Cells.Item(2, 9).Value = "Loket 3 (Wiji)"
$ws.Cells.Item(2, 11).Value = "Ereg/NPWP"

# Row 3 (NO=2)
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).NumberFormat = "@"
$ws.Cells.Item(3, 2).Value = "2022-11-16"
$ws.Cells.Item(3, 3).Value = "yesaya"
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "617694500712000"
$ws.Cells.Item(3, 5).Value = "Orang Pribadi"
$ws.Cells.Item(3, 6).Value = "jl . tangkuhis no. 30"
$ws.Cells.Item(3, 7).Value = "karyawan swasta"
$ws.Cells.Item(3, 9).Value = "Loket 3 (Wiji)"
$ws.Cells.Item(3, 11).Value = "Ereg/NPWP"

# Row 4 (NO=3)
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).NumberFormat = "@"
$ws.Cells.Item(4, 2).Value = "2022-11-16"
$ws.Cells.Item(4, 3).Value = "INDRA GANDI"
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "087179222712000"
$ws.Cells.Item(4, 5).Value = "Orang Pribadi"
$ws.Cells.Item(4, 6).Value = "KARUING 001 001, KAB. KATINGAN, KARUING, KAMIPANG, KALIMANTAN TENGAH"
$ws.Cells.Item(4, 7).Value = "KLU WP Unknown"
$ws.Cells.Item(4, 9).Value = "Loket 3 (Wiji)"
$ws.Cells.Item(4, 11).Value = "Ereg/NPWP,SPT Tahunan / SPT Masa"

# Row 5 (NO=4)
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).NumberFormat = "@"
$ws.Cells.Item(5, 2).Value = "2022-11-16"
$ws.Cells.Item(5, 3).Value = "PUTRI NOR MENTARY"
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "762306215712000"
$ws.Cells.Item(5, 5).Value = "Orang Pribadi"
$ws.Cells.Item(5, 6).Value = "DESA PETAK BAHANDANG RT 001 RW 001, KAB. KATINGAN, PETAK BAHANDANG, TASIK PAYAWAN, KALIMANTAN TENGAH"
$ws.Cells.Item(5, 7).Value = "PEGAWAI SWASTA"
$ws.Cells.Item(5, 9).Value = "Loket 3 (Wiji)"
$ws.Cells.Item(5, 11).Value = "Ereg/NPWP,SPT Tahunan / SPT Masa"

# Row 6 (NO=5)
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).NumberFormat = "@"
$ws.Cells.Item(6, 2).Value = "2022-11-15"
$ws.Cells.Item(6, 3).Value = "NI KADE ANGGARAINI"
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "086515517712000"
$ws.Cells.Item(6, 5).Value = "Orang Pribadi"
$ws.Cells.Item(6, 6).Value = "JALAN KALI SAMBA NO. 025 RT 002, KAB. KATINGAN, TUMBANG KAMAN, SANAMAN MANTIKEI, KALIMANTAN TENGAH"
$ws.Cells.Item(6, 7).Value = "KLU WP Unknown"
$ws.Cells.Item(6, 9).Value = "Loket 3 (Wiji)"
$ws.Cells.Item(6, 11).Value = "Ereg/NPWP"

# Row 7 (NO=6)
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).NumberFormat = "@"
$ws.Cells.Item(7, 2).Value = "2022-11-15"
$ws.Cells.Item(7, 3).Value = "KRISTIN"
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "253991624712000"
$ws.Cells.Item(7, 5).Value = "Orang Pribadi"
$ws.Cells.Item(7, 6).Value = "JL. A. YANI RT/RW 007/03, KAB. KATINGAN, SAMBA KAHAYAN, KATINGAN TENGAH, KALIMANTAN TENGAH"
$ws.Cells.Item(7, 7).Value = "KLU WP Unknown"
$ws.Cells.Item(7, 8).Value = "krist712@mailnesia.com"
$ws.Cells.Item(7, 9).Value = "Loket 4 (Azriel)"
$ws.Cells.Item(7, 11).Value = "Ereg/NPWP"
$ws.Cells.Item(7, 12).Value = "aktivasi efin / djp online & cetak ulang"

# Row 8 (NO=7)
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).NumberFormat = "@"
$ws.Cells.Item(8, 2).Value = "2022-11-15"
$ws.Cells.Item(8, 3).Value = "HANLI"
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "617407275712000"
$ws.Cells.Item(8, 5).Value = "Orang Pribadi"
$ws.Cells.Item(8, 6).Value = "TUMBANG PANGGO"
$ws.Cells.Item(8, 7).Value = "KEBUN, "
$ws.Cells.Item(8, 8).Value = "hanli712@yopmail.com"
$ws.Cells.Item(8, 9).Value = "Loket 4 (Azriel)"
$ws.Cells.Item(8, 10).Value = "NPWP elektronik dikirim lewat WA. WP akan jemput (jika ada waktu)"
$ws.Cells.Item(8, 11).Value = "Ereg/NPWP"
$ws.Cells.Item(8, 12).Value = "Daftar NPWP"

# Row 9 (NO=8)
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).NumberFormat = "@"
$ws.Cells.Item(9, 2).Value = "2022-11-15"
$ws.Cells.Item(9, 3).Value = "RAHMAD FAUZI"
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "844602268712000"
$ws.Cells.Item(9, 5).Value = "Orang Pribadi"
$ws.Cells.Item(9, 6).Value = "JL. TJILIK RIWUT KM. 14"
$ws.Cells.Item(9, 7).Value = "Dagang"
$ws.Cells.Item(9, 8).Value = "rahmad712@mailnesia.com"
$ws.Cells.Item(9, 9).Value = "Loket 4 (Azriel)"
$ws.Cells.Item(9, 11).Value = "Ereg/NPWP"
$ws.Cells.Item(9, 12).Value = "cetak ulang npwp (hilang)"

# Row 10 (NO=9)
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).NumberFormat = "@"
$ws.Cells.Item(10, 2).Value = "2022-11-15"
$ws.Cells.Item(10, 3).Value = "WANDRIE"
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "321696098712000"
$ws.Cells.Item(10, 5).Value = "Orang Pribadi"
$ws.Cells.Item(10, 6).Value = "DESA TEWANG KARANGAN NO. 01 RT 001 RW 001 KEL. PULAU MALAN KASUNGAN 74453 KEC. PULAU MALAN, KAB. KATINGAN, TEWANG KARANGAN, PULAU MALAN, KALIMANTAN TENGAH"
$ws.Cells.Item(10, 7).Value = "KLU WP Unknown"
$ws.Cells.Item(10, 8).Value = "wandrie712@yopmail.com"
$ws.Cells.Item(10, 9).Value = "Loket 4 (Azriel)"
$ws.Cells.Item(10, 10).Value = "KK : 6206041111100168"
$ws.Cells.Item(10, 11).Value = "Ereg/NPWP"
$ws.Cells.Item(10, 12).Value = "Daftar NPWP"

# Row 11 (NO=10)
$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).NumberFormat = "@"
$ws.Cells.Item(11, 2).Value = "2022-11-14"
$ws.Cells.Item(11, 3).Value = "Wiji Thukul"
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "666666666712000"
$ws.Cells.Item(11, 5).Value = "Orang Pribadi"
$ws.Cells.Item(11, 6).Value = "Jl Kasongan Lama"
$ws.Cells.Item(11, 7).Value = "Tukang bubur"
$ws.Cells.Item(11, 8).Value = "bruh@bruhmail.com"
$ws.Cells.Item(11, 9).Value = "Loket 3 (Wiji)"
$ws.Cells.Item(11, 10).Value = "wp diminta untuk pulang ke rahmatulloh"
$ws.Cells.Item(11, 11).Value = "E-Billing,SPT Tahunan / SPT Masa"

# Renumber the pushed-down original rows (now rows 12-15): NO 1-4 -> 11-14
$ws.Cells.Item(12, 1).Value = 11
$ws.Cells.Item(13, 1).Value = 12
$ws.Cells.Item(14, 1).Value = 13
$ws.Cells.Item(15, 1).Value = 14

Write-Host "Edit complete"
